# "misc download bug fixed"
# The author's name shown on Sheet1 (cell A2, driven by the shared string
# table) was wrong - it gets corrected here. After the edit the user's
# cursor/selection moved on to the next cell (A3), which is the normal
# Excel behaviour after committing a cell edit with Enter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Stefan hawking"

$ws.Range("A3").Select()
